# "Generate Report for Handback" -- mark the localization status workbook
# as handed back (in sync with en-US) and record the handback file name /
# timestamp for the zh-cn and de-de target sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status text for both locales ---------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("I2").Value = "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.md"
$zhcn.Range("J2").Value = "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.04104383c7affb9ef1f9a4f05e6882cb10d5b276.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-13 09:11:38"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/fdcd04a489f5c3620c8420509f597f2c31982caa/e2e/2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.md", [Type]::Missing, "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.md") | Out-Null
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet ---------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("I2").Value = "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.md"
$dede.Range("J2").Value = "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.04104383c7affb9ef1f9a4f05e6882cb10d5b276.de-de.xlf"
$dede.Range("K2").Value = "2016-08-13 09:11:47"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/fdcd04a489f5c3620c8420509f597f2c31982caa/e2e/2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.md", [Type]::Missing, "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3.md") | Out-Null
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

Write-Output "Handback report generated."
